$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.842.80"
$ws.Range("E2").Value = "  +0.32%  "

$ws.Range("D3").Value = "2.570.27"
$ws.Range("E3").Value = "  +1.63%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.56"
$ws.Range("E5").Value = "  -0.95%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.27"
$ws.Range("E6").Value = "  +2.88%  "

$ws.Range("E7").Value = "  +0.16%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.533"
$ws.Range("E9").Value = "  +0.32%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.72"
$ws.Range("E10").Value = "  -0.13%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0809"
$ws.Range("E11").Value = "  +0.60%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.44"
$ws.Range("E12").Value = "  -1.05%  "

$ws.Range("D13").Value = "2.956.86"
$ws.Range("E13").Value = "  +1.29%  "

$ws.Range("E14").Value = "  -1.34%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.83"
$ws.Range("E15").Value = "  +4.82%  "

$ws.Range("D16").Value = "2.618.15"
$ws.Range("E16").Value = "  +2.59%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.842"
$ws.Range("E17").Value = "  -0.48%  "

$ws.Range("D18").Value = "42.858.47"
$ws.Range("E18").Value = "  +0.09%  "

$ws.Range("E19").Value = "  -1.70%  "

$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.44"
$ws.Range("E20").Value = "  -2.77%  "

$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "0.0₃0959"
$ws.Range("E21").Value = "  -0.15%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.74"
$ws.Range("E22").Value = "  +0.18%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "248.64"
$ws.Range("E23").Value = "  -1.01%  "

$ws.Range("E24").Value = "  -0.23%  "

$ws.Range("E25").Value = "  -0.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "27.05"
$ws.Range("E26").Value = "  +2.27%  "

$ws.Range("E27").Value = "  +0.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.41"
$ws.Range("E28").Value = "  -0.10%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.82"
$ws.Range("E29").Value = "  -1.28%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.22"
$ws.Range("E30").Value = "  -1.56%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.57"
$ws.Range("E31").Value = "  +1.99%  "

$ws.Range("E32").Value = "  -1.97%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0798"
$ws.Range("E33").Value = "  +3.02%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.11"
$ws.Range("E34").Value = "  -2.09%  "

$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.68"
$ws.Range("E35").Value = "  -0.32%  "

$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.30"
$ws.Range("E36").Value = "  -0.75%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.65"
$ws.Range("E37").Value = "  -0.85%  "

$ws.Range("E38").Value = "  +10.87%  "

$ws.Range("E39").Value = "  -0.12%  "

$ws.Range("E40").Value = "  -0.07%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.83"
$ws.Range("E41").Value = "  +2.55%  "

$ws.Range("E42").Value = "  +7.87%  "

$ws.Range("E43").Value = "  -0.11%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0302"
$ws.Range("E44").Value = "  -0.27%  "

$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.21"
$ws.Range("E45").Value = "  -0.81%  "

$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "1.996.50"
$ws.Range("E46").Value = "  -1.62%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.00"
$ws.Range("E47").Value = "  -0.51%  "

$ws.Range("D48").Value = "2.809.09"
$ws.Range("E48").Value = "  +1.27%  "

$ws.Range("E49").Value = "  +2.71%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "81.56"
$ws.Range("E50").Value = "  -3.44%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "74.23"
$ws.Range("E51").Value = "  -0.56%  "
